$d = $word.ActiveDocument

# First paragraph of the document (the **ID__...__ID** marker paragraph)
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with a 5pt "space" (distance
# from text) value, matching an added <w:pBdr> with w:space="5" on each edge.
$borders = $p1.Range.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Change the left indent from 120 twips (=6pt) to 225 twips (=11.25pt).
$p1.Format.LeftIndent = 11.25

# Replace the marker text and drop the trailing space run: the paragraph
# used to contain two runs - "**ID__AFFARS_5325_topic_12__ID**" followed by
# a run with a single space " ". We search (scoped to this paragraph) for
# the text plus trailing space and replace it with the new marker text (no
# trailing space), which merges/collapses the two runs into a single run
# with the desired text.
$find = $p1.Range.Find
$find.ClearFormatting()
$null = $find.Execute("**ID__AFFARS_5325_topic_12__ID** ", $true, $false, $false, $false, $false, `
               $true, 1, $false, "**ID__AFFARS_5325_1001__ID**", 2)
